$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "135  p p p"

# Clear out old content first (A2:A3 previously had data)
$ws.Cells.Clear()

# Row 2
$ws.Range("A2").Value = "Тест1"
$ws.Range("B2").Value = "пор"

# Row 3 (order chosen to reproduce the shared-strings table order from the diff)
$ws.Range("D3").Value = "Україна"
$ws.Range("A3").Value = "Осв пр 1"
$ws.Range("B3").Value = "kz"
$ws.Range("C3").Value = "rti"
$ws.Range("H3").Value = "dddd"

# Update the active selection
$ws.Range("E5").Select()
